$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stats")

$ws.Range("D2").Value = 0.00313084339722991
$ws.Range("E2").Value = 1.439415222033858
$ws.Range("G2").Value = 0.09222432691603899
$ws.Range("H2").Value = 0.7298780623823404
$ws.Range("I2").Value = 0.2350507462397218
$ws.Range("J2").Value = 0.2684049438685179
$ws.Range("K2").Value = 0.02796277310699224
$ws.Range("D3").Value = 0.06545410817489028
$ws.Range("E3").Value = 0.3947789939120412
$ws.Range("G3").Value = 0.0104007925838232
$ws.Range("H3").Value = 0.1424757433123887
$ws.Range("I3").Value = 0.02403827849775553
$ws.Range("J3").Value = 0.2072149789892137
$ws.Range("K3").Value = 0.002761213108897209
$ws.Range("D4").Value = 0.06641591154038906
$ws.Range("E4").Value = 0.3943201447837055
$ws.Range("G4").Value = 0.0105892731808126
$ws.Range("H4").Value = 0.1423162538558245
$ws.Range("I4").Value = 0.02404126944020391
$ws.Range("J4").Value = 0.2065588254481554
$ws.Range("K4").Value = 0.002788963727653027
$ws.Range("D5").Value = 0.004563175607472658
$ws.Range("E5").Value = 1.440237573347986
$ws.Range("G5").Value = 0.09064242616295815
$ws.Range("H5").Value = 0.72658471763134
$ws.Range("I5").Value = 0.2432540557347238
$ws.Range("J5").Value = 0.2673658370040357
$ws.Range("K5").Value = 0.02728631859645247
$ws.Range("D6").Value = 0.06115334015339613
$ws.Range("E6").Value = 0.5087484568357468
$ws.Range("G6").Value = 0.01396965002641082
$ws.Range("H6").Value = 0.1617975165136158
$ws.Range("I6").Value = 0.1007672557607293
$ws.Range("J6").Value = 0.2169510298408568
$ws.Range("K6").Value = 0.004315282683819532
$ws.Range("E7").Value = 26.88509621098638
$ws.Range("D8").Value = 0.002693851012736559
$ws.Range("E8").Value = 1.247489012777805
$ws.Range("G8").Value = 0.07787654967978597
$ws.Range("H8").Value = 0.6121293497271836
$ws.Range("I8").Value = 0.2155747562646866
$ws.Range("J8").Value = 0.2461630539037287
$ws.Range("K8").Value = 0.02345868293195963
$ws.Range("D9").Value = 0.06393868941813707
$ws.Range("E9").Value = 0.4722560863010585
$ws.Range("G9").Value = 0.0120175639167428
$ws.Range("H9").Value = 0.1602854132652283
$ws.Range("I9").Value = 0.02254338981583714
$ws.Range("J9").Value = 0.2635422684252262
$ws.Range("K9").Value = 0.003504962660372257
$ws.Range("D10").Value = 0.06312387529760599
$ws.Range("E10").Value = 0.4153491077013314
$ws.Range("G10").Value = 0.01154327439144254
$ws.Range("H10").Value = 0.1537308827973902
$ws.Range("I10").Value = 0.02199699822813272
$ws.Range("J10").Value = 0.2147464766167104
$ws.Range("K10").Value = 0.003328687977045774
$ws.Range("D11").Value = 0.00381794385612011
$ws.Range("E11").Value = 1.27689443808049
$ws.Range("G11").Value = 0.07650809921324253
$ws.Range("H11").Value = 0.6029225932434201
$ws.Range("I11").Value = 0.2206326364539564
$ws.Range("J11").Value = 0.2818358112126589
$ws.Range("K11").Value = 0.0234526083804667
$ws.Range("D12").Value = 0.07000776287168264
$ws.Range("E12").Value = 0.5625158869661391
$ws.Range("G12").Value = 0.01577644562348723
$ws.Range("H12").Value = 0.182832152582705
$ws.Range("I12").Value = 0.1092394990846515
$ws.Range("J12").Value = 0.2370333350263536
$ws.Range("K12").Value = 0.004977410659193993
$ws.Range("E13").Value = 26.13077858230099
$ws.Range("D14").Value = 0.002912347204983234
$ws.Range("E14").Value = 1.343452117405832
$ws.Range("G14").Value = 0.08505043829791248
$ws.Range("H14").Value = 0.671003706054762
$ws.Range("I14").Value = 0.2253127512522042
$ws.Range("J14").Value = 0.2572839988861233
$ws.Range("K14").Value = 0.02571072801947594
$ws.Range("D15").Value = 0.06469639879651368
$ws.Range("E15").Value = 0.4335175401065499
$ws.Range("G15").Value = 0.011209178250283
$ws.Range("H15").Value = 0.1513805782888085
$ws.Range("I15").Value = 0.02329083415679634
$ws.Range("J15").Value = 0.23537862370722
$ws.Range("K15").Value = 0.003133087884634733
$ws.Range("D16").Value = 0.06476989341899753
$ws.Range("E16").Value = 0.4048346262425184
$ws.Range("G16").Value = 0.01106627378612757
$ws.Range("H16").Value = 0.1480235683266073
$ws.Range("I16").Value = 0.02301913383416831
$ws.Range("J16").Value = 0.2106526510324329
$ws.Range("K16").Value = 0.003058825852349401
$ws.Range("D17").Value = 0.004190559731796384
$ws.Range("E17").Value = 1.358566005714238
$ws.Range("G17").Value = 0.08357526268810034
$ws.Range("H17").Value = 0.6647536554373801
$ws.Range("I17").Value = 0.2319433460943401
$ws.Range("J17").Value = 0.2746008241083473
$ws.Range("K17").Value = 0.02536946348845959
$ws.Range("D18").Value = 0.06558055151253939
$ws.Range("E18").Value = 0.5356321719009429
$ws.Range("G18").Value = 0.01487304782494903
$ws.Range("H18").Value = 0.1723148345481604
$ws.Range("I18").Value = 0.1050033774226904
$ws.Range("J18").Value = 0.2269921824336052
$ws.Range("K18").Value = 0.004646346671506763
$ws.Range("E19").Value = 26.50793739664368
